$d = $word.ActiveDocument

# --- 1. Fix the heading: merge "ECM1410 Cover " + "page" + " " into one run ---
# Replace the whole heading paragraph's range (spanning all 3 runs + the
# proofErr markers) with the single consolidated phrase - this collapses
# everything down to one run and drops the proofErr bookmarks.
$headingPara = $d.Paragraphs.Item(1)
$headingRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End)
$headingRange.Text = "ECM1410 Cover page "

# --- 2. Fill in the last (previously empty) row of the second table ---
$tbl = $d.Tables.Item(2)
$lastRow = $tbl.Rows.Count
$tbl.Cell($lastRow, 1).Range.Text = "03/3/23"
$tbl.Cell($lastRow, 2).Range.Text = "16:00"
$tbl.Cell($lastRow, 3).Range.Text = "20mins"
$tbl.Cell($lastRow, 4).Range.Text = "Observer"
$tbl.Cell($lastRow, 5).Range.Text = "Driver"
